$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "Subnational indicators explorer "
$ws.Range("B3").Value = "Compare a local authority and the UK average (median) local authority by different indicators, such as weekly pay and healthy life expectancy."
$ws.Range("C3").Value = "ONS"
$ws.Range("D3").Value = "https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/subnationalindicatorsexplorer/2022-01-06"
$ws.Range("E3").Value = "Publicly available"

# Row 4
$ws.Range("A4").Value = "Association of Colleges (AoC)"
$ws.Range("B4").Value = "Lists useful data sources related to colleges"
$ws.Range("C4").Value = "AoC"
$ws.Range("D4").Value = "https://www.aoc.co.uk/research-unit/data-sources"
$ws.Range("E4").Value = "Publicly available"

# Row 5
$ws.Range("A5").Value = "Census"
$ws.Range("B5").Value = "Contains a range of topics including labour market and education by local authority. "
$ws.Range("C5").Value = "ONS"
$ws.Range("D5").Value = "https://census.gov.uk/local-authorities"
$ws.Range("E5").Value = "Publicly available"

# Row 6
$ws.Range("A6").Value = "NOMIS"
$ws.Range("B6").Value = "Query labour market data at national, regional and local levels"
$ws.Range("C6").Value = "ONS"
$ws.Range("D6").Value = "https://www.nomisweb.co.uk/"
$ws.Range("E6").Value = "Publicly available, with more functionality if you create an account. "

# Row 7
$ws.Range("A7").Value = "Explore Education Statistics"
$ws.Range("B7").Value = "A range of educational data published by the Department of Education "
$ws.Range("C7").Value = "DfE"
$ws.Range("D7").Value = "https://explore-education-statistics.service.gov.uk/"
$ws.Range("E7").Value = "Publicly available "

$ws.Range("A:E").AutoFit() | Out-Null

$ws.Range("F20").Select() | Out-Null
